$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the modified timestamp
$ws.Range("B20").Value = "2022-06-20T11:36:40+00:00"

# Fill in the skos:altLabel (column C) values for the new variable rows
$ws.Range("C24").Value = "PM_RT"
$ws.Range("C25").Value = "LDT_RT"
$ws.Range("C26").Value = "OG_RT"
$ws.Range("C28").Value = "PM_accuracy"
$ws.Range("C29").Value = "LDT_accuracy"
$ws.Range("C30").Value = "OG_accuracy"
